# Edit applies the diff changes: re-positions several annotation shapes on
# slide 12 (SlideID 279) to better align with an updated scatter-plot image,
# splits the bold "Remark:" sentence out of the slide-12 title textbox into
# its own smaller, non-bold callout textbox, and nudges a few oval/picture
# annotations on slides 14, 15 and 17 (SlideIDs 283, 281, 282).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 12 (SlideID 279)
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)

# Shape id=4 ("כותרת 1" / title textbox) - move up slightly and drop the
# trailing "Remark: ..." paragraph (it becomes its own textbox below).
$shTitle = $s12.Shapes.Item(3)
$shTitle.Top = 1697232 / 12700
$tr = $shTitle.TextFrame.TextRange
$tr.Text = "Route B: " + [char]13 + "Analyzing the entire filtered Properties" + [char]8217 + " population"

# Picture 2 (id=1026) - reposition.
$shPic = $s12.Shapes.Item(4)
$shPic.Left = 2724839 / 12700
$shPic.Top = 2431739 / 12700

# Straight Arrow Connector 6 (id=7) - reposition.
$shCxn1 = $s12.Shapes.Item(5)
$shCxn1.Left = 4561490 / 12700
$shCxn1.Top = 5574959 / 12700

# Oval 8 (id=9) - reposition.
$shOval1 = $s12.Shapes.Item(6)
$shOval1.Left = 2724839 / 12700
$shOval1.Top = 3530354 / 12700

# Oval 9 (id=10) - reposition + resize.
$shOval2 = $s12.Shapes.Item(7)
$shOval2.Left = 9297217 / 12700
$shOval2.Top = 5263557 / 12700
$shOval2.Width = 260164 / 12700
$shOval2.Height = 311402 / 12700

# TextBox 12 (id=13, "Outliers") - reposition.
$shTb12 = $s12.Shapes.Item(8)
$shTb12.Left = 8400237 / 12700
$shTb12.Top = 4659646 / 12700

# Straight Arrow Connector 14 (id=15) - reposition + resize.
$shCxn2 = $s12.Shapes.Item(9)
$shCxn2.Left = 9002487 / 12700
$shCxn2.Top = 4994253 / 12700
$shCxn2.Width = 294730 / 12700
$shCxn2.Height = 269304 / 12700

# New textbox carrying the "Remark: ..." sentence split out of the title
# shape above. The slide-local shape id/name counter increments from 2, so
# add (and discard) a throwaway textbox first to land on id=3 / "TextBox 2"
# exactly like the recorded edit.
$throwaway = $s12.Shapes.AddTextbox(1, 0, 0, 10, 10)
$throwaway.Delete()

$remarkBox = $s12.Shapes.AddTextbox(1, 1382108 / 12700, 6000219 / 12700, 6096000 / 12700, 369332 / 12700)
$remarkBox.Fill.Visible = $false
$remarkBox.TextFrame.WordWrap = $true
$remarkBox.TextFrame.AutoSize = 1

$remarkTr = $remarkBox.TextFrame.TextRange
$remarkTr.Text = "Remark: The depended variable on the x-axis for convenience "
$remarkTr.ParagraphFormat.Alignment = 1
$remarkTr.Font.Size = 18
$remarkTr.Font.NameComplexScript = "+mn-cs"

$remarkLabel = $remarkTr.Characters(1, 6)
$remarkLabel.Font.Underline = $true

# ---------------------------------------------------------------------
# Slide 14 (SlideID 283)
# ---------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$shOvalS14 = $s14.Shapes.Item(6)
$shOvalS14.Width = 248963 / 12700
$shOvalS14.Height = 1062447 / 12700

# ---------------------------------------------------------------------
# Slide 15 (SlideID 281)
# ---------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$shOvalS15 = $s15.Shapes.Item(6)
$shOvalS15.Top = 3679706 / 12700

# ---------------------------------------------------------------------
# Slide 17 (SlideID 282)
# ---------------------------------------------------------------------
$s17 = $p.Slides.Item(17)
$shPicS17 = $s17.Shapes.Item(4)
$shPicS17.Left = 1545393 / 12700

$shOvalS17 = $s17.Shapes.Item(6)
$shOvalS17.Top = 3438350 / 12700
